# Update data: 4 September 2020
# Appends the newest month (2020-08-01, Excel serial 44044) of unemployment
# data to both the "Canada" (sheet1) and "Province" (sheet2) worksheets.

$wb = $excel.ActiveWorkbook

$dateFormat = "d-mmm-yy"
$newDate = 44044

# ---------------------------------------------------------------------------
# Sheet "Canada": append one row (row 9) with the national figures.
# ---------------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A9").Value = $newDate
$wsCanada.Range("A9").NumberFormat = $dateFormat
$wsCanada.Range("B9").Value = "Canada"
$wsCanada.Range("B9").NumberFormat = $dateFormat
$wsCanada.Range("C9").Value = 76.5
$wsCanada.Range("D9").Value = 2046.9

$wsCanada.Range("C10").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Province": append ten rows (rows 72-81), one per province, for the
# same reporting month, following the existing layout (Newfoundland &
# Labrador first - with the date style also applied to column B - then the
# remaining provinces in their established order).
# ---------------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 72; Name = "Newfoundland & Labrador"; C = -3.6;  D = 32.1;   DateStyleOnB = $true  },
    @{ Row = 73; Name = "Prince Edward Island";     C = 19.7;  D = 9.1;    DateStyleOnB = $false },
    @{ Row = 74; Name = "Nova Scotia";               C = 27.3;  D = 50.8;   DateStyleOnB = $false },
    @{ Row = 75; Name = "New Brunswick";             C = 9;     D = 36.2;   DateStyleOnB = $false },
    @{ Row = 76; Name = "Quebec";                    C = 83.6;  D = 398.4;  DateStyleOnB = $false },
    @{ Row = 77; Name = "Ontario";                   C = 90;    D = 841.4;  DateStyleOnB = $false },
    @{ Row = 78; Name = "Manitoba";                  C = 45.8;  D = 56;     DateStyleOnB = $false },
    @{ Row = 79; Name = "Saskatchewan";               C = 53.9;  D = 47.7;   DateStyleOnB = $false },
    @{ Row = 80; Name = "Alberta";                    C = 59.2;  D = 290.4;  DateStyleOnB = $false },
    @{ Row = 81; Name = "British Columbia";           C = 112.5; D = 284.9;  DateStyleOnB = $false }
)

foreach ($r in $provinceRows) {
    $rowNum = $r.Row

    $wsProvince.Range("A$rowNum").Value = $newDate
    $wsProvince.Range("A$rowNum").NumberFormat = $dateFormat

    $wsProvince.Range("B$rowNum").Value = $r.Name
    if ($r.DateStyleOnB) {
        $wsProvince.Range("B$rowNum").NumberFormat = $dateFormat
    }

    $wsProvince.Range("C$rowNum").Value = $r.C
    $wsProvince.Range("D$rowNum").Value = $r.D
}

$wsProvince.Activate() | Out-Null
$wsProvince.Range("C82").Select() | Out-Null
